# Commit: "Fruta / hortaliza, semanal"
# Adds one new week of price data (two quality-grade rows: "Primera" and
# "Segunda") for Albahaca at Vega Central Mapocho de Santiago, inserted
# right before the existing row 118 — pushing all the subsequent rows
# down by two and growing the sheet from A1:R183 to A1:R185.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new records by inserting two blank rows at 118.
$ws.Rows("118:119").Insert()

# New row 118 — "Primera" grade.
$ws.Range("A118").Value = 9
$ws.Range("B118").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C118").Value = "Metropolitana"
$ws.Range("D118").Value = 44455
$ws.Range("E118").Value = 13
$ws.Range("F118").Value = 100112052
$ws.Range("G118").Value = "Albahaca"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 178
$ws.Range("K118").Value = 5500
$ws.Range("L118").Value = 6000
$ws.Range("M118").Value = 5750
$ws.Range("N118").Value = "$/paquete"
$ws.Range("O118").Value = "Región de Arica y Parinacota"
$ws.Range("P118").Value = 5750
$ws.Range("Q118").Value = 1
$ws.Range("R118").Value = "Hortaliza"

# New row 119 — "Segunda" grade.
$ws.Range("A119").Value = 9
$ws.Range("B119").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C119").Value = "Metropolitana"
$ws.Range("D119").Value = 44455
$ws.Range("E119").Value = 13
$ws.Range("F119").Value = 100112052
$ws.Range("G119").Value = "Albahaca"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Segunda"
$ws.Range("J119").Value = 97
$ws.Range("K119").Value = 4000
$ws.Range("L119").Value = 4500
$ws.Range("M119").Value = 4247
$ws.Range("N119").Value = "$/paquete"
$ws.Range("O119").Value = "Región de Arica y Parinacota"
$ws.Range("P119").Value = 4247
$ws.Range("Q119").Value = 1
$ws.Range("R119").Value = "Hortaliza"
